$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new teacher row (row 5): name in A5, email hyperlink in B5 -
# mirrors the existing rows (2-4) which each hold a teacher name and a
# mailto hyperlink to idriselbasaur@gmail.com.
$ws.Range("A5").Value = "Mr. Sam"
$ws.Range("B5").Value = "idriselbasaur@gmail.com"

$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:idriselbasaur@gmail.com", "", "", "idriselbasaur@gmail.com")

# Adding a hyperlink registers Excel's built-in "Hyperlink" cell style;
# bring B5's look back in line with the other email cells (B2:B4) and
# drop the now-unused named style so formatting stays consistent.
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$wb.Styles("Hyperlink").Delete()

$ws.Range("B5").Select()
